# Fixed falsed negative transitions: cannot go to recovered or dead straight
# from this state.
#
# The "FN,NQ" row (row 14) and "FN,Q" row (row 15) previously had entries
# for transitions into the "RA,NQ"/"RA,Q" (recovered) and "D" (dead) states.
# Those transitions are invalid, so the corresponding cell contents are
# cleared while keeping their existing formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14 ("FN,NQ"): clear Q14 (omega^r) and S14 (omega^D)
$ws.Range("Q14").ClearContents()
$ws.Range("S14").ClearContents()

# Row 15 ("FN,Q"): clear R15 (omega^r) and S15 (omega^D)
$ws.Range("R15").ClearContents()
$ws.Range("S15").ClearContents()

# Update the active selection to match the saved view state
$ws.Range("S15").Select()
